# Price update for 2026-02-07 — append the latest scraped row
# (Date, Price, Discount, Incredible) to the bottom of the tracking sheet.

$wb = $excel.ActiveWorkbook
$wsName = $wb.ActiveSheet.Name

# Stage the new values on a throwaway sheet first and force them to be
# plain text (not auto-converted to a date serial / number by Excel),
# matching how every other row in the sheet is stored (t="s" shared
# strings with no special number formatting).
$tmp = $wb.Worksheets.Add()
$tmp.Name = "__tmp_stage__"
$staging = $tmp.Range("A1:D1")
$staging.NumberFormat = "@"
$tmp.Range("A1").Value = "2026-02-07"
$tmp.Range("B1").Value = "9999000"
$tmp.Range("C1").Value = "28"
$tmp.Range("D1").Value = "0"

$staging.Copy()

$ws = $wb.Worksheets.Item($wsName)
$lastRow = $ws.Cells(1048576, 1).End(-4162).Row
$targetRow = $lastRow + 1
$ws.Cells.Item($targetRow, 1).PasteSpecial(-4163)
$excel.CutCopyMode = $false

$tmp.Delete() | Out-Null
